# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Albahaca" (Vega Modelo de Temuco) above the
# existing row 232, shifting the following rows down by one (old row 232
# becomes 233, ..., old row 244 becomes 245).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 232; existing rows 232:244 shift down to 233:245.
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with this week's record.
$ws.Range("A232").Value = 10
$ws.Range("B232").Value = "Vega Modelo de Temuco"
$ws.Range("C232").Value = "La Araucanía"
$ws.Range("D232").Value = 44753
$ws.Range("E232").Value = 9
$ws.Range("F232").Value = 100112052
$ws.Range("G232").Value = "Albahaca"
$ws.Range("H232").Value = "Sin especificar"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 80
$ws.Range("K232").Value = 6000
$ws.Range("L232").Value = 6000
$ws.Range("M232").Value = 6000
$ws.Range("N232").Value = "`$/paquete"
$ws.Range("O232").Value = "Región de Arica y Parinacota"
$ws.Range("P232").Value = 6000
$ws.Range("Q232").Value = 1
$ws.Range("R232").Value = "Hortaliza"
